# Atualização do fluxo de potência e função da bateria. Modularização do Main
#
# BESS sheet: rework the battery parameters.
#   Einit        -> Einit(%)          (value unchanged)
#   Emax         -> Cmax              (value unchanged)
#   Emin(%)      -> SOC_min(%)        (value 20 -> 10)
#   (new column) -> SOC_max(%) = 90   (inserted before the old Efficiency column)
#   Efficiency stays last (shifts one column to the right automatically)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BESS")

# Insert a new column before the old "Efficiency" column (col I / 9) so the
# existing Efficiency header+value shift right to column J, picking up
# their formatting as-is.
$ws.Columns.Item(9).Insert()

# Rename the existing headers that keep their position.
$ws.Range("F1").Value = "Einit(%)"

# Update the existing numeric value that changed.
$ws.Range("H1").Value = "SOC_min(%)"
$ws.Range("H2").Value = 10

# Fill in the newly inserted SOC_max(%) column. The inserted column already
# picked up H's formatting (bold header, right-aligned value), so a plain
# value assignment is enough to match the look of SOC_min(%).
$ws.Range("I1").Value = "SOC_max(%)"
$ws.Range("I2").Value = 90

$ws.Range("G1").Value = "Cmax"

# Resize the two columns to fit their new contents.
$ws.Columns.Item(8).ColumnWidth = 12.61
$ws.Columns.Item(9).ColumnWidth = 13.17

# BESS becomes the active sheet / tab, with H3 selected (Generator, which
# used to be active, loses tabSelected automatically).
$ws.Activate()
$ws.Range("H3").Select()
